$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 2601.8684239554491
$ws.Range("B1").Value = 2133.4903394187268
$ws.Range("C1").Value = 2059.7275084469711
$ws.Range("A2").Value = 2283.417581285787
$ws.Range("B2").Value = 1494.9049736678382
$ws.Range("C2").Value = 1700.7225660525737
$ws.Range("A3").Value = 2875.1540281488974
$ws.Range("B3").Value = 2193.7678485077363
$ws.Range("C3").Value = 2138.889820529434
$ws.Range("A4").Value = 2716.5390107578428
$ws.Range("B4").Value = 2202.8785414699923
$ws.Range("C4").Value = 2210.8070497552576
$ws.Range("A5").Value = 2921.2417263228745
$ws.Range("B5").Value = 2147.9433594646393
$ws.Range("C5").Value = 2198.802912162032
$ws.Range("A6").Value = 2781.5017078405867
$ws.Range("B6").Value = 2303.5786862321388
$ws.Range("C6").Value = 2387.5638409348853
$ws.Range("A7").Value = 2446.1630423940978
$ws.Range("B7").Value = 2357.4685080142717
$ws.Range("C7").Value = 2273.8995891981413
$ws.Range("A8").Value = 2884.6357911765899
$ws.Range("B8").Value = 2506.2688726822835
$ws.Range("C8").Value = 2322.0383854150141
$ws.Range("A9").Value = 3159.6791931575567
$ws.Range("B9").Value = 2252.6717554136508
$ws.Range("C9").Value = 2035.8814965623105
$ws.Range("A10").Value = 2442.5674573863539
$ws.Range("B10").Value = 1443.5584288180928
$ws.Range("C10").Value = 1577.3614355783966
$ws.Range("A11").Value = 2163.7771558412483
$ws.Range("B11").Value = 1666.0900062610754
$ws.Range("C11").Value = 1596.1313243617089
$ws.Range("A12").Value = 3066.1447316021095
$ws.Range("B12").Value = 2504.8774066351243
$ws.Range("C12").Value = 2362.1266634315061
$ws.Range("A13").Value = 3145.3317022853021
$ws.Range("B13").Value = 2517.3950040905615
$ws.Range("C13").Value = 2258.3215857169948
$ws.Range("A14").Value = 3159.0007489776867
$ws.Range("B14").Value = 2565.8378342902743
$ws.Range("C14").Value = 2289.3734640242396
$ws.Range("A15").Value = 3127.5956997100775
$ws.Range("B15").Value = 2515.7552776066377
$ws.Range("C15").Value = 2389.3838300551997
$ws.Range("A16").Value = 3015.4375493997914
$ws.Range("B16").Value = 2161.4677253310683
$ws.Range("C16").Value = 1838.1350059045205
